# Auto update stock data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [new date, new EBITDA value (or $null if unchanged)]
$updates = @(
    @{ Row = 2;  Date = "2025/12/25"; Ebitda = "6.67" },
    @{ Row = 8;  Date = "2025/12/25"; Ebitda = "8.56" },
    @{ Row = 14; Date = "2025/12/25"; Ebitda = $null },
    @{ Row = 20; Date = "2025/12/25"; Ebitda = $null },
    @{ Row = 26; Date = "2025/12/25"; Ebitda = "11.26" },
    @{ Row = 32; Date = "2025/12/25"; Ebitda = "27.82" },
    @{ Row = 38; Date = "2025/12/25"; Ebitda = $null },
    @{ Row = 44; Date = "2025/12/25"; Ebitda = "11.25" },
    @{ Row = 50; Date = "2025/12/25"; Ebitda = "11.60" },
    @{ Row = 56; Date = "2025/12/25"; Ebitda = "32.38" },
    @{ Row = 62; Date = "2025/12/25"; Ebitda = "11.69" },
    @{ Row = 68; Date = "2025/12/25"; Ebitda = "13.21" },
    @{ Row = 74; Date = "2025/12/25"; Ebitda = "16.70" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $u.Date
    if ($null -ne $u.Ebitda) {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $u.Ebitda
    }
}
